$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text so values like
# "11.50" or "0.00001070" keep their exact formatting instead of
# being coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.105.26"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").Value = "1.806.21"
$ws.Range("E3").Value = "  +2.72%  "

$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "324.20"
$ws.Range("E5").Value = "  -0.27%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "0.4292"
$ws.Range("E7").Value = "  -4.00%  "

$ws.Range("D8").Value = "0.3638"
$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").Value = "44.86"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("D10").Value = "0.07583"
$ws.Range("E10").Value = "  +0.39%  "

$ws.Range("D11").Value = "1.131"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").Value = "21.68"
$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "6.253"
$ws.Range("E14").Value = "  +0.64%  "

$ws.Range("D15").Value = "7.406"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").Value = "1.829.52"
$ws.Range("E16").Value = "  +3.76%  "

$ws.Range("D17").Value = "93.37"
$ws.Range("E17").Value = "  +6.16%  "

$ws.Range("D18").Value = "0.00001070"
$ws.Range("E18").Value = "  -0.43%  "

$ws.Range("D19").Value = "0.06386"
$ws.Range("E19").Value = "  +2.20%  "

$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "17.31"
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "6.174"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").Value = "28.163.25"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D24").Value = "11.50"
$ws.Range("E24").Value = "  -1.49%  "

$ws.Range("D25").Value = "2.144"
$ws.Range("E25").Value = "  -7.60%  "

$ws.Range("D26").Value = "160.40"
$ws.Range("E26").Value = "  +4.93%  "

$ws.Range("D27").Value = "20.53"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("D28").Value = "2.031.66"
$ws.Range("E28").Value = "  +3.66%  "

$ws.Range("D29").Value = "2.242"
$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("D30").Value = "130.02"
$ws.Range("E30").Value = "  +1.12%  "

$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  -3.45%  "

$ws.Range("D32").Value = "5.935"
$ws.Range("E32").Value = "  +3.23%  "

$ws.Range("D33").Value = "0.09075"
$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("D34").Value = "3.538"
$ws.Range("E34").Value = "  -2.87%  "

$ws.Range("D35").Value = "12.84"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("D36").Value = "0.02385"
$ws.Range("E36").Value = "  +2.34%  "

$ws.Range("D37").Value = "5.153"
$ws.Range("E37").Value = "  +1.39%  "

$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "0.6513"
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2137"
$ws.Range("E39").Value = "  -1.73%  "

$ws.Range("D40").Value = "0.06140"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("D41").Value = "1.199"
$ws.Range("E41").Value = "  -0.07%  "

$ws.Range("E42").Value = "  +0.33%  "

$ws.Range("D43").Value = "7.980"
$ws.Range("E43").Value = "  -0.32%  "

$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "13.77"
$ws.Range("E45").Value = "  -0.35%  "

$ws.Range("D46").Value = "0.6037"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("D47").Value = "3.722"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "125.74"
$ws.Range("E48").Value = "  -0.52%  "

$ws.Range("D49").Value = "1.996"
$ws.Range("E49").Value = "  +0.26%  "

$ws.Range("D50").Value = "1.175"
$ws.Range("E50").Value = "  +3.62%  "

$ws.Range("D51").Value = "0.06992"
$ws.Range("E51").Value = "  +1.26%  "
